# Add new columns I (I0) and J (IF) to the sheet, mirroring the style of the
# existing header cells (e.g. H1) and filling in the per-row numeric data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header's formatting (bold, centered, bordered) onto the
# two new header cells so they reuse the same cell style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-13 for columns I and J
$data = @{
    2  = @(1, 6)
    3  = @(1, 5)
    4  = @(1, 4)
    5  = @(1, 3)
    6  = @(1, 6)
    7  = @(7, 8)
    8  = @(5, 5)
    9  = @(8, 9)
    10 = @(1, 3)
    11 = @(9, 9)
    12 = @(6, 8)
    13 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}

Write-Output "Added I0/IF columns"
